$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) figures.
# Values are written as text (matching the inline-string cells already
# used in this sheet) rather than being auto-converted to numbers/percentages.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "261.39"
Set-TextValue "E2" "0.73%"
Set-TextValue "D3" "27.12"
Set-TextValue "E3" "0.71%"
Set-TextValue "D4" "4.724"
Set-TextValue "E4" "0.61%"
Set-TextValue "D5" "0.06202"
Set-TextValue "E5" "2.52%"
Set-TextValue "D6" "6.726"
Set-TextValue "E6" "0.61%"
Set-TextValue "D7" "0.8518"
Set-TextValue "E7" "-0.95%"
Set-TextValue "D8" "0.9062"
Set-TextValue "E8" "-1.92%"
Set-TextValue "E9" "0.82%"
Set-TextValue "D10" "0.04925"
Set-TextValue "E10" "-4.22%"
Set-TextValue "D11" "0.07098"
Set-TextValue "E11" "0.23%"
Set-TextValue "D12" "0.03166"
Set-TextValue "E12" "3.20%"
Set-TextValue "D13" "0.09059"
Set-TextValue "E13" "-0.84%"
Set-TextValue "D14" "0.001528"
Set-TextValue "E14" "-0.37%"
Set-TextValue "D15" "0.0006142"
Set-TextValue "E15" "1.35%"
Set-TextValue "D16" "0.006140"
Set-TextValue "E16" "1.88%"
Set-TextValue "D17" "3.466"
Set-TextValue "E17" "-0.36%"
Set-TextValue "D18" "3.168"
Set-TextValue "E18" "-0.02%"
Set-TextValue "E19" "-0.34%"
Set-TextValue "E20" "-0.68%"
Set-TextValue "E21" "-1.27%"
Set-TextValue "E22" "-0.12%"
Set-TextValue "E23" "-0.34%"
Set-TextValue "E24" "0.17%"
Set-TextValue "E25" "2.51%"
Set-TextValue "E26" "0.18%"
Set-TextValue "D40" "0.03913"
Set-TextValue "E40" "1.45%"
Set-TextValue "D41" "0.1113"
Set-TextValue "E41" "-0.19%"
Set-TextValue "D42" "0.004138"
Set-TextValue "E42" "2.31%"
Set-TextValue "E43" "-0.65%"
Set-TextValue "E44" "-8.87%"
Set-TextValue "D45" "0.00005174"
Set-TextValue "E45" "-0.42%"
Set-TextValue "D47" "0.03591"
Set-TextValue "E47" "-34.12%"
Set-TextValue "D48" "0.1699"
Set-TextValue "E48" "25.54%"
